$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

function Add-HandbackColumns {
    param($ws, $mdCellSrc, $xlfCellSrc, $targetCell, $handbackCell, $statusCell, $datetimeCell, $datetimeValue, $mdUrl, $mdDisplay, $xlfUrl, $xlfDisplay)

    # Update the status to reflect that the file has been handed back and is in sync.
    $ws.Range($statusCell).Value = $statusHandedBack

    # Populate "Latest Target File" and "Latest Handback File" columns with hyperlinks
    # pointing at the same source/handoff files (mirroring columns A and C).
    $ws.Hyperlinks.Add($ws.Range($targetCell), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range($handbackCell), $xlfUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $xlfDisplay)

    # Record the handback datetime.
    $ws.Range($datetimeCell).Value = $datetimeValue
}

# ----- Overview sheet: status column also reflects the handback -----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

Add-HandbackColumns $wsZh "A2" "C2" "E2" "F2" "B2" "G2" "2016-03-11 00:13:28" "https://github.com/OpenLocalizationTest/oltest/blob/9062a1160a00fe267c812f0cfad2318063d5fa8d/e2e/71788b6e-7678-40a3-a49b-60fde762a7cb.md" "71788b6e-7678-40a3-a49b-60fde762a7cb.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71b498191cb660c746d18b0ad86b6b5343380c2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/71788b6e-7678-40a3-a49b-60fde762a7cb.2fa3ecf2447f2631aa2405b1c1b059808f615f75.zh-cn.xlf" "71788b6e-7678-40a3-a49b-60fde762a7cb.2fa3ecf2447f2631aa2405b1c1b059808f615f75.zh-cn.xlf"

Add-HandbackColumns $wsZh "A3" "C3" "E3" "F3" "B3" "G3" "2016-03-11 00:13:28" "https://github.com/OpenLocalizationTest/oltest/blob/9062a1160a00fe267c812f0cfad2318063d5fa8d/e2e/c09f4f10-9639-48c8-8abf-ad89a4a12872.md" "c09f4f10-9639-48c8-8abf-ad89a4a12872.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/71b498191cb660c746d18b0ad86b6b5343380c2d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/c09f4f10-9639-48c8-8abf-ad89a4a12872.1dfcbfe505d40e2d1f6854ae528577f895df6d5a.zh-cn.xlf" "c09f4f10-9639-48c8-8abf-ad89a4a12872.1dfcbfe505d40e2d1f6854ae528577f895df6d5a.zh-cn.xlf"

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

Add-HandbackColumns $wsDe "A2" "C2" "E2" "F2" "B2" "G2" "2016-03-11 00:13:47" "https://github.com/OpenLocalizationTest/oltest/blob/9062a1160a00fe267c812f0cfad2318063d5fa8d/e2e/71788b6e-7678-40a3-a49b-60fde762a7cb.md" "71788b6e-7678-40a3-a49b-60fde762a7cb.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e06bd8935739b7e975dd1dbc6a4678013bb9964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/71788b6e-7678-40a3-a49b-60fde762a7cb.2fa3ecf2447f2631aa2405b1c1b059808f615f75.de-de.xlf" "71788b6e-7678-40a3-a49b-60fde762a7cb.2fa3ecf2447f2631aa2405b1c1b059808f615f75.de-de.xlf"

Add-HandbackColumns $wsDe "A3" "C3" "E3" "F3" "B3" "G3" "2016-03-11 00:13:47" "https://github.com/OpenLocalizationTest/oltest/blob/9062a1160a00fe267c812f0cfad2318063d5fa8d/e2e/c09f4f10-9639-48c8-8abf-ad89a4a12872.md" "c09f4f10-9639-48c8-8abf-ad89a4a12872.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3e06bd8935739b7e975dd1dbc6a4678013bb9964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/c09f4f10-9639-48c8-8abf-ad89a4a12872.1dfcbfe505d40e2d1f6854ae528577f895df6d5a.de-de.xlf" "c09f4f10-9639-48c8-8abf-ad89a4a12872.1dfcbfe505d40e2d1f6854ae528577f895df6d5a.de-de.xlf"
